# Updated code to use ECI frame for thrust vector
# Adds three new columns (L, M, N) holding the ECI X/Y/Z Thrust Component
# values alongside the existing RIC thrust component columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (row 1)
$ws.Range("L1").Value = "ECI X Thrust Component"
$ws.Range("M1").Value = "ECI Y Thrust Component"
$ws.Range("N1").Value = "ECI Z Thrust Component"

# New ECI thrust component data for each burn segment row
$eciData = @{
    4  = @(-0.6561, 0.50032, -0.56495)
    6  = @(-0.7392, 0.66354, 0.11472)
    9  = @(0.10832, 0.1314, -0.985387)
    12 = @(0.11428, 0.12571, -0.98546)
    15 = @(0.1194, 0.11965, -0.9856)
    18 = @(0.12439, 0.11327, -0.9857)
    21 = @(0.12944, 0.10666, -0.98583)
    24 = @(0.1338, 0.09972, -0.98597)
    27 = @(0.138, 0.092556, -0.98609)
    30 = @(0.14173, 0.085123, -0.986237)
}

foreach ($row in $eciData.Keys) {
    $values = $eciData[$row]
    $ws.Cells.Item($row, 12).Value = $values[0]
    $ws.Cells.Item($row, 13).Value = $values[1]
    $ws.Cells.Item($row, 14).Value = $values[2]
}

# Update the active view/selection to match the edited region
try {
    $excel.ActiveWindow.ScrollRow = 14
    $excel.ActiveWindow.ScrollColumn = 2
} catch {
    # Scroll-position automation isn't exposed everywhere; ignore if missing.
}
$ws.Range("N31").Select()
